$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1111126
$ws.Range("I6").Value = 1111126
$ws.Range("K6").Value = 3333378
$ws.Range("M6").Value = -3333266
$ws.Range("H58").Value = 10664.182
$ws.Range("I58").Value = 1539.125
$ws.Range("J58").Value = 34997.668
$ws.Range("K58").Value = 4617.375
$ws.Range("L58").Value = 104993.004
$ws.Range("M58").Value = -4467.375
$ws.Range("N58").Value = -105293.004
$ws.Range("H80").Value = 706.4167
$ws.Range("I80").Value = 127
$ws.Range("J80").Value = 996.125
$ws.Range("K80").Value = 381
$ws.Range("L80").Value = 2988.375
$ws.Range("M80").Value = 617
$ws.Range("N80").Value = -4984.375
$ws.Range("H83").Value = 706.4167
$ws.Range("I83").Value = 127
$ws.Range("J83").Value = 996.125
$ws.Range("K83").Value = 1143
$ws.Range("L83").Value = 8965.125
$ws.Range("M83").Value = 3849
$ws.Range("N83").Value = -18949.125
$ws.Range("H86").Value = 4612
$ws.Range("J86").Value = 5334
$ws.Range("L86").Value = 5334
$ws.Range("N86").Value = -7580
$ws.Range("H89").Value = 4612
$ws.Range("J89").Value = 5334
$ws.Range("L89").Value = 26670
$ws.Range("N89").Value = -37902
$ws.Range("H100").Value = 1480.8462
$ws.Range("I100").Value = 1643.0952
$ws.Range("J100").Value = 799.4
$ws.Range("K100").Value = 1643.0952
$ws.Range("L100").Value = 799.4
$ws.Range("M100").Value = -1102.0952
$ws.Range("N100").Value = -1881.4
$ws.Range("H112").Value = 38642.9
$ws.Range("I112").Value = 2473.5
$ws.Range("J112").Value = 41322.11
$ws.Range("K112").Value = 7420.5
$ws.Range("L112").Value = 123966.33
$ws.Range("M112").Value = -6312.5
$ws.Range("N112").Value = -126182.33
$ws.Range("H118").Value = 1631.1666
$ws.Range("I118").Value = 508.44446
$ws.Range("K118").Value = 1525.33338
$ws.Range("M118").Value = 131.66662
$ws.Range("H127").Value = 1232.0526
$ws.Range("I127").Value = 846.26666
$ws.Range("J127").Value = 2678.75
$ws.Range("K127").Value = 2538.79998
$ws.Range("L127").Value = 8036.25
$ws.Range("M127").Value = 2421.20002
$ws.Range("N127").Value = -17956.25
$ws.Range("H131").Value = 9257.087
$ws.Range("I131").Value = 1828.9445
$ws.Range("J131").Value = 35998.4
$ws.Range("K131").Value = 5486.833500000001
$ws.Range("L131").Value = 107995.2
$ws.Range("M131").Value = -446.8335000000006
$ws.Range("N131").Value = -118075.2
$ws.Range("H132").Value = 1473.475
$ws.Range("I132").Value = 1434.3334
$ws.Range("K132").Value = 4303.0002
$ws.Range("M132").Value = -1773.0002
$ws.Range("H138").Value = 7249358
$ws.Range("I138").Value = 1357.8889
$ws.Range("J138").Value = 9807476
$ws.Range("K138").Value = 4073.6667
$ws.Range("L138").Value = 29422428
$ws.Range("M138").Value = 1066.3333
$ws.Range("N138").Value = -29432708
$ws.Range("H141").Value = 966.2308
$ws.Range("I141").Value = 996.4545000000001
$ws.Range("J141").Value = 800
$ws.Range("K141").Value = 2989.3635
$ws.Range("L141").Value = 2400
$ws.Range("M141").Value = 2190.6365
$ws.Range("N141").Value = -12760

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 11425
$ws.Range("I74").Value = 2043.2858
$ws.Range("K74").Value = 2043.2858
$ws.Range("M74").Value = -1169.2858
$ws.Range("H77").Value = 11425
$ws.Range("I77").Value = 2043.2858
$ws.Range("K77").Value = 10216.429
$ws.Range("M77").Value = -5848.429
$ws.Range("H97").Value = 42958.582
$ws.Range("I97").Value = 749
$ws.Range("J97").Value = 113307.89
$ws.Range("K97").Value = 749
$ws.Range("L97").Value = 113307.89
$ws.Range("M97").Value = -253
$ws.Range("N97").Value = -114299.89
$ws.Range("H102").Value = 2704.3684
$ws.Range("I102").Value = 2691.5334
$ws.Range("J102").Value = 2752.5
$ws.Range("K102").Value = 2691.5334
$ws.Range("L102").Value = 2752.5
$ws.Range("M102").Value = -1069.5334
$ws.Range("N102").Value = -5996.5
$ws.Range("H132").Value = 3082.3333
$ws.Range("I132").Value = 2472.7368
$ws.Range("K132").Value = 7418.2104
$ws.Range("M132").Value = -4888.2104

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 420.33334
$ws.Range("J5").Value = 99.666664
$ws.Range("L5").Value = 99.666664
$ws.Range("N5").Value = -325.666664
$ws.Range("H94").Value = 1161.7333
$ws.Range("I94").Value = 991.6
$ws.Range("J94").Value = 1502
$ws.Range("K94").Value = 991.6
$ws.Range("L94").Value = 1502
$ws.Range("M94").Value = -540.6
$ws.Range("N94").Value = -2404
$ws.Range("H99").Value = 91119.30499999999
$ws.Range("I99").Value = 60809.117
$ws.Range("K99").Value = 60809.117
$ws.Range("M99").Value = -59311.117
$ws.Range("H126").Value = 90779
$ws.Range("J126").Value = 90779
$ws.Range("L126").Value = 90779
$ws.Range("N126").Value = -100659
$ws.Range("H134").Value = 2259.9517
$ws.Range("I134").Value = 2185.3
$ws.Range("K134").Value = 6555.900000000001
$ws.Range("M134").Value = -4020.900000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 4779.8
$ws.Range("I12").Value = 2474.75
$ws.Range("K12").Value = 2474.75
$ws.Range("M12").Value = -2304.75
$ws.Range("H60").Value = 17979.4
$ws.Range("I60").Value = 16633
$ws.Range("J60").Value = 19999
$ws.Range("K60").Value = 16633
$ws.Range("L60").Value = 19999
$ws.Range("M60").Value = -16122
$ws.Range("N60").Value = -21021
$ws.Range("H68").Value = 99999.5
$ws.Range("J68").Value = 99999.5
$ws.Range("L68").Value = 99999.5
$ws.Range("N68").Value = -101497.5
$ws.Range("H71").Value = 99999.5
$ws.Range("J71").Value = 99999.5
$ws.Range("L71").Value = 299998.5
$ws.Range("N71").Value = -307486.5
$ws.Range("H81").Value = 40000
$ws.Range("J81").Value = 40000
$ws.Range("L81").Value = 40000
$ws.Range("N81").Value = -41996
$ws.Range("H84").Value = 40000
$ws.Range("J84").Value = 40000
$ws.Range("L84").Value = 120000
$ws.Range("N84").Value = -129984
$ws.Range("H86").Value = 8114.8887
$ws.Range("J86").Value = 7574.25
$ws.Range("L86").Value = 7574.25
$ws.Range("N86").Value = -9820.25
$ws.Range("H89").Value = 8114.8887
$ws.Range("J89").Value = 7574.25
$ws.Range("L89").Value = 37871.25
$ws.Range("N89").Value = -49103.25
$ws.Range("H105").Value = 1795.8334
$ws.Range("I105").Value = 1795.8334
$ws.Range("K105").Value = 1795.8334
$ws.Range("M105").Value = -48.83339999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 49052584
$ws.Range("I4").Value = 57166880
$ws.Range("J4").Value = 366783.34
$ws.Range("K4").Value = 171500640
$ws.Range("L4").Value = 1100350.02
$ws.Range("M4").Value = -171500528
$ws.Range("N4").Value = -1100574.02
$ws.Range("H23").Value = 568.5263
$ws.Range("J23").Value = 625.38464
$ws.Range("L23").Value = 1876.15392
$ws.Range("N23").Value = -2346.15392
$ws.Range("H122").Value = 1129.421
$ws.Range("I122").Value = 419.8
$ws.Range("K122").Value = 3778.2
$ws.Range("M122").Value = -1328.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 27778484
$ws.Range("I102").Value = 718.2258
$ws.Range("K102").Value = 718.2258
$ws.Range("M102").Value = 903.7742
$ws.Range("H132").Value = 3921.6553
$ws.Range("I132").Value = 3465.8696
$ws.Range("K132").Value = 10397.6088
$ws.Range("M132").Value = -7867.6088
$ws.Range("H141").Value = 67309.836
$ws.Range("J141").Value = 67309.836
$ws.Range("L141").Value = 67309.836
$ws.Range("N141").Value = -77669.836

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 50030
$ws.Range("I38").Value = 50030
$ws.Range("K38").Value = 50030
$ws.Range("M38").Value = -49620
$ws.Range("H93").Value = 1857.5
$ws.Range("I93").Value = 1833.5
$ws.Range("J93").Value = 2001.5
$ws.Range("K93").Value = 1833.5
$ws.Range("L93").Value = 2001.5
$ws.Range("M93").Value = -585.5
$ws.Range("N93").Value = -4497.5
$ws.Range("H132").Value = 3144.28
$ws.Range("I132").Value = 2947.7896
$ws.Range("J132").Value = 3766.5
$ws.Range("K132").Value = 8843.3688
$ws.Range("L132").Value = 11299.5
$ws.Range("M132").Value = -6313.3688
$ws.Range("N132").Value = -16359.5
$ws.Range("H136").Value = 4692.643
$ws.Range("I136").Value = 3563.3635
$ws.Range("K136").Value = 10690.0905
$ws.Range("M136").Value = -8140.0905

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3031.5
$ws.Range("I96").Value = 1851.5
$ws.Range("J96").Value = 3621.5
$ws.Range("K96").Value = 1851.5
$ws.Range("L96").Value = 3621.5
$ws.Range("M96").Value = -478.5
$ws.Range("N96").Value = -6367.5
$ws.Range("H100").Value = 800.6799999999999
$ws.Range("I100").Value = 782.8333
$ws.Range("J100").Value = 846.5714
$ws.Range("K100").Value = 1565.6666
$ws.Range("L100").Value = 1693.1428
$ws.Range("M100").Value = -1024.6666
$ws.Range("N100").Value = -2775.1428
$ws.Range("H107").Value = 253000
$ws.Range("I107").Value = 3000
$ws.Range("K107").Value = 9000
$ws.Range("M107").Value = -7080
$ws.Range("H136").Value = 1869
$ws.Range("I136").Value = 1721.92
$ws.Range("K136").Value = 5165.76
$ws.Range("M136").Value = -2615.76
